$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the shared string used by the old D2 value into the new I1 header text,
# and clear out the old D2/B2/C2 cell contents.
$ws.Range("D2").Value = $null
$ws.Range("B2").Value = $null
$ws.Range("C2").Value = $null

# Add the new "EpiLength" column header.
$ws.Range("I1").Value = "EpiLength"

# Fill column A with the sequential row numbers for the generalInfo table.
$ws.Range("A2").Value = 1
$ws.Range("A3").Value = 2
$ws.Range("A4").Value = 3
$ws.Range("A5").Value = 4
$ws.Range("A6").Value = 5
$ws.Range("A7").Value = 6
$ws.Range("A8").Value = 7

# Update the active selection like in the edited workbook.
$ws.Range("D14").Select()
